$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D-column cells whose new value would otherwise be
# auto-converted to a Number by Excel (losing exact text representation,
# e.g. trailing zeros). Applied only to the cells that need it.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '53.758.75'
$ws.Range("E2").Value = '  -8.25%  '
$ws.Range("D3").Value = '2.417.09'
$ws.Range("E3").Value = '  -13.75%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '462.58'
$ws.Range("E5").Value = '  -6.81%  '
$ws.Range("D6").Value = '131.06'
$ws.Range("E6").Value = '  -1.24%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").Value = '0.484'
$ws.Range("E8").Value = '  -7.68%  '
$ws.Range("D9").Value = '2.430.76'
$ws.Range("E9").Value = '  -13.30%  '
$ws.Range("D10").Value = '0.0948'
$ws.Range("E10").Value = '  -6.14%  '
$ws.Range("D11").Value = '5.30'
$ws.Range("E11").Value = '  -9.69%  '
$ws.Range("D12").Value = '0.315'
$ws.Range("E12").Value = '  -8.03%  '
$ws.Range("D13").Value = '0.121'
$ws.Range("E13").Value = '  -4.24%  '
$ws.Range("D14").Value = '2.858.25'
$ws.Range("E14").Value = '  -13.44%  '
$ws.Range("D15").Value = '53.675.07'
$ws.Range("E15").Value = '  -8.64%  '
$ws.Range("D16").Value = '0.0000133'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '19.66'
$ws.Range("E17").Value = '  -7.35%  '
$ws.Range("D18").Value = '2.457.51'
$ws.Range("E18").Value = '  -12.52%  '
$ws.Range("D19").Value = '4.17'
$ws.Range("E19").Value = '  -10.03%  '
$ws.Range("D20").Value = '310.62'
$ws.Range("E20").Value = '  -9.35%  '
$ws.Range("D21").Value = '9.36'
$ws.Range("E21").Value = '  -13.78%  '
$ws.Range("D22").Value = '0.992'
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").Value = '5.71'
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("D24").Value = '5.35'
$ws.Range("E24").Value = '  -13.11%  '
$ws.Range("D25").Value = '56.18'
$ws.Range("E25").Value = '  -10.22%  '
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").Value = '2.559.07'
$ws.Range("E27").Value = '  -13.06%  '
$ws.Range("D28").Value = '0.381'
$ws.Range("E28").Value = '  -9.42%  '
$ws.Range("D29").Value = '0.152'
$ws.Range("E29").Value = '  -10.06%  '
$ws.Range("D30").Value = '7.17'
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").Value = '0.0₃0720'
$ws.Range("E32").Value = '  -7.72%  '
$ws.Range("D33").Value = '150.36'
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").Value = '17.59'
$ws.Range("E34").Value = '  -6.42%  '
$ws.Range("E35").Value = '  -10.03%  '
$ws.Range("D36").Value = '4.99'
$ws.Range("E36").Value = '  -4.64%  '
$ws.Range("D37").Value = '3.52'
$ws.Range("E37").Value = '  -13.11%  '
$ws.Range("D38").Value = '1.05'
$ws.Range("E38").Value = '  -4.80%  '
$ws.Range("D39").Value = '0.788'
$ws.Range("E39").Value = '  -10.87%  '
$ws.Range("D40").Value = '33.43'
$ws.Range("E40").Value = '  -9.20%  '
$ws.Range("D41").Value = '0.993'
$ws.Range("E41").Value = '  -0.74%  '
$ws.Range("D42").Value = '0.596'
$ws.Range("E42").Value = '  -4.18%  '
$ws.Range("D43").Value = '0.0526'
$ws.Range("E43").Value = '  -4.49%  '
$ws.Range("D44").Value = '3.25'
$ws.Range("E44").Value = '  -6.01%  '
$ws.Range("D45").Value = '10.19'
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("D46").Value = '1.23'
$ws.Range("E46").Value = '  -6.61%  '
$ws.Range("D47").Value = '1.950.96'
$ws.Range("E47").Value = '  -11.50%  '
$ws.Range("D48").Value = '0.0219'
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("D49").Value = '0.0867'
$ws.Range("E49").Value = '  -1.68%  '
$ws.Range("D50").Value = '4.24'
$ws.Range("E50").Value = '  -6.07%  '
$ws.Range("D51").Value = '16.58'
$ws.Range("E51").Value = '  -12.00%  '
